# Continuando con la memoria...
# Adds a new hyperparameter-search result row (row 29) to the main table and
# a small "semilla" (seed) lookup table in columns H:I (rows 19-24) on the
# single worksheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking string as literal TEXT (shared string),
# matching what Excel does when such text is pasted in rather than typed
# (keeps cell type "s" and avoids creating a throw-away cell style, since a
# plain `.Value = "0.8259"` assignment gets auto-coerced to a number by the
# host, same as typing it into the UI would).
function Set-TextValue {
    param($cell, [string]$text)

    $scratch = $ws.Range("ZZ1")
    $escaped = $text.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
    $scratch.ClearContents() | Out-Null
}

# --- New row 29 in the main results table -----------------------------
Set-TextValue $ws.Range("E29") "700 iter"
Set-TextValue $ws.Range("F29") "0.119214"
Set-TextValue $ws.Range("A29") "0.8253"
Set-TextValue $ws.Range("B29") "0.3"
$ws.Range("C29").Value = 15
Set-TextValue $ws.Range("D29") "0.02"

# --- New "semilla" lookup table in H19:I24 -----------------------------
Set-TextValue $ws.Range("H19") "semilla"
Set-TextValue $ws.Range("I19") "Test"

$ws.Range("H20").Value = 123
Set-TextValue $ws.Range("I20") "0.8259"

$ws.Range("H21").Value = 1234
Set-TextValue $ws.Range("I21") "0.8260"

$ws.Range("H22").Value = 12345
Set-TextValue $ws.Range("I22") "0.8253"

$ws.Range("H23").Value = 123456
Set-TextValue $ws.Range("I23") "0.8257"

$ws.Range("H24").Value = 1244
Set-TextValue $ws.Range("I24") "0.8242"

# Match the saved selection/active cell recorded in the workbook.
$ws.Range("I25").Select() | Out-Null
